$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (column C) bumps from 45245 to 45246 for every existing data row
for ($r = 2; $r -le 26; $r++) {
    $ws.Cells.Item($r, 3).Value = 45246
}

# The last logged item ("A 57074-2023") has dropped out of the source feed,
# so its row is removed entirely.
$ws.Rows.Item(27).Delete()

# Row 26 becomes the new last row and reverts to the default (non-custom) row
# height, matching how freshly-appended rows are left unformatted.
$ws.Rows.Item(26).AutoFit()
